$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 21
$ws.Range("B2").Value = 6

$ws.Range("A3").Value = 22
$ws.Range("B3").Value = 5

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 1
